$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $savedStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $savedStyle
}

Set-TextValue "D2" "69.583.22"
Set-TextValue "E2" "  +1.73%  "
Set-TextValue "D3" "2.442.31"
Set-TextValue "E3" "  +0.14%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.15%  "
Set-TextValue "D5" "564.54"
Set-TextValue "E5" "  +1.21%  "
Set-TextValue "D6" "165.79"
Set-TextValue "E6" "  +1.94%  "
Set-TextValue "E7" "  -0.11%  "
Set-TextValue "D8" "0.510"
Set-TextValue "E8" "  -0.18%  "
Set-TextValue "D9" "0.175"
Set-TextValue "E9" "  +11.37%  "
Set-TextValue "E10" "  -1.46%  "
Set-TextValue "E11" "  +2.05%  "
Set-TextValue "D12" "4.66"
Set-TextValue "E12" "  -3.02%  "
Set-TextValue "D13" "0.0000182"
Set-TextValue "E13" "  +7.13%  "
Set-TextValue "D14" "69.401.42"
Set-TextValue "E14" "  +1.57%  "
Set-TextValue "D15" "2.882.17"
Set-TextValue "E15" "  -0.81%  "
Set-TextValue "D16" "23.88"
Set-TextValue "E16" "  +2.58%  "
Set-TextValue "D17" "2.444.94"
Set-TextValue "E17" "  +0.01%  "
Set-TextValue "E18" "  +3.22%  "
Set-TextValue "D19" "340.64"
Set-TextValue "E19" "  +1.14%  "
Set-TextValue "D20" "7.08"
Set-TextValue "E20" "  +2.80%  "
Set-TextValue "D21" "3.87"
Set-TextValue "E21" "  +2.05%  "
Set-TextValue "D22" "2.01"
Set-TextValue "E22" "  +6.48%  "
Set-TextValue "E23" "  -0.03%  "
Set-TextValue "D24" "66.16"
Set-TextValue "E24" "  -0.89%  "
Set-TextValue "E25" "  +4.96%  "
Set-TextValue "D26" "2.577.51"
Set-TextValue "E26" "  +0.44%  "
Set-TextValue "D27" "8.45"
Set-TextValue "E27" "  +3.58%  "
Set-TextValue "D28" "0.958"
Set-TextValue "E28" "  -4.09%  "
Set-TextValue "D29" "0.0₃0848"
Set-TextValue "E29" "  +4.05%  "
Set-TextValue "D30" "7.28"
Set-TextValue "E30" "  +1.34%  "
Set-TextValue "D31" "1.25"
Set-TextValue "E31" "  +9.36%  "
Set-TextValue "D32" "445.94"
Set-TextValue "E32" "  +4.88%  "
Set-TextValue "D33" "0.999"
Set-TextValue "E33" "  -0.10%  "
Set-TextValue "E34" "  +0.44%  "
Set-TextValue "D35" "161.30"
Set-TextValue "E35" "  +0.46%  "
Set-TextValue "D36" "19.05"
Set-TextValue "E36" "  +0.33%  "
Set-TextValue "E38" "  +3.72%  "
Set-TextValue "D39" "18.07"
Set-TextValue "E39" "  +1.75%  "
Set-TextValue "D40" "0.306"
Set-TextValue "E40" "  +3.15%  "
Set-TextValue "B41" "Stacks"
Set-TextValue "C41" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D41" "1.54"
Set-TextValue "E41" "  +4.66%  "
Set-TextValue "B42" "RenderToken"
Set-TextValue "C42" "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D42" "4.44"
Set-TextValue "E42" "  +1.84%  "
Set-TextValue "D43" "1.09"
Set-TextValue "E43" "  +2.19%  "
Set-TextValue "D44" "2.13"
Set-TextValue "E44" "  +5.87%  "
Set-TextValue "E45" "  +0.92%  "
Set-TextValue "D46" "131.44"
Set-TextValue "E46" "  +1.46%  "
Set-TextValue "D47" "0.0723"
Set-TextValue "E47" "  +1.00%  "
Set-TextValue "D48" "0.487"
Set-TextValue "E48" "  +1.34%  "
Set-TextValue "E49" "  -0.11%  "
Set-TextValue "D50" "0.0928"
Set-TextValue "E50" "  +1.21%  "
Set-TextValue "E51" "  +2.69%  "
